$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7 currently lists "Grand Piano" with a hyperlink to the iconshock
# musical-instruments icon page. We rename that row's icon-name to
# "Kettle Drum" (it now describes the kettle-drum icon, reusing the same
# source link) and add a new row 8 underneath for "Grand Piano", which
# keeps pointing at the very same iconshock page.

$iconshockUrl = $ws.Range("B7").Text

# Row 7: Grand Piano -> Kettle Drum (link/style/URL text unchanged).
$ws.Range("A7").Value = "Kettle Drum"

# Row 8 (new): Grand Piano, with its own hyperlink to the iconshock page.
$ws.Range("A8").Value = "Grand Piano"
$ws.Range("B8").Value = $iconshockUrl
$ws.Hyperlinks.Add($ws.Range("B8"), $iconshockUrl) | Out-Null
$ws.Range("B8").Style = "Hyperlink"
